$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values in column A (rows 2-201)
$ws.Range("A2").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 3
$ws.Range("A7").Value = 3
$ws.Range("A8").Value = 3
$ws.Range("A9").Value = 2
$ws.Range("A11").Value = 3
$ws.Range("A13").Value = 3
$ws.Range("A14").Value = 2
$ws.Range("A15").Value = 2
$ws.Range("A16").Value = 1
$ws.Range("A17").Value = 1
$ws.Range("A19").Value = 1
$ws.Range("A20").Value = 3
$ws.Range("A21").Value = 3
$ws.Range("A22").Value = 3
$ws.Range("A23").Value = 3
$ws.Range("A28").Value = 3
$ws.Range("A29").Value = 3
$ws.Range("A30").Value = 3
$ws.Range("A32").Value = 3
$ws.Range("A33").Value = 2
$ws.Range("A35").Value = 3
$ws.Range("A36").Value = 1
$ws.Range("A37").Value = 3
$ws.Range("A39").Value = 3
$ws.Range("A43").Value = 1
$ws.Range("A44").Value = 1
$ws.Range("A45").Value = 2
$ws.Range("A47").Value = 3
$ws.Range("A50").Value = 1
$ws.Range("A52").Value = 3
$ws.Range("A54").Value = 1
$ws.Range("A55").Value = 3
$ws.Range("A56").Value = 3
$ws.Range("A57").Value = 2
$ws.Range("A58").Value = 2
$ws.Range("A61").Value = 1
$ws.Range("A62").Value = 3
$ws.Range("A63").Value = 1
$ws.Range("A64").Value = 1
$ws.Range("A66").Value = 3
$ws.Range("A69").Value = 2
$ws.Range("A71").Value = 3
$ws.Range("A72").Value = 2
$ws.Range("A73").Value = 3
$ws.Range("A74").Value = 1
$ws.Range("A75").Value = 1
$ws.Range("A76").Value = 1
$ws.Range("A78").Value = 2
$ws.Range("A79").Value = 3
$ws.Range("A80").Value = 2
$ws.Range("A82").Value = 1
$ws.Range("A85").Value = 3
$ws.Range("A86").Value = 3
$ws.Range("A88").Value = 2
$ws.Range("A89").Value = 1
$ws.Range("A91").Value = 3
$ws.Range("A94").Value = 3
$ws.Range("A96").Value = 2
$ws.Range("A97").Value = 1
$ws.Range("A99").Value = 3
$ws.Range("A100").Value = 1
$ws.Range("A101").Value = 2
$ws.Range("A102").Value = 2
$ws.Range("A103").Value = 3
$ws.Range("A104").Value = 1
$ws.Range("A106").Value = 2
$ws.Range("A109").Value = 3
$ws.Range("A110").Value = 3
$ws.Range("A111").Value = 3
$ws.Range("A112").Value = 1
$ws.Range("A113").Value = 3
$ws.Range("A114").Value = 1
$ws.Range("A115").Value = 3
$ws.Range("A117").Value = 2
$ws.Range("A118").Value = 2
$ws.Range("A120").Value = 3
$ws.Range("A126").Value = 3
$ws.Range("A127").Value = 3
$ws.Range("A130").Value = 1
$ws.Range("A131").Value = 2
$ws.Range("A132").Value = 2
$ws.Range("A133").Value = 3
$ws.Range("A137").Value = 1
$ws.Range("A138").Value = 1
$ws.Range("A139").Value = 2
$ws.Range("A140").Value = 3
$ws.Range("A141").Value = 2
$ws.Range("A142").Value = 1
$ws.Range("A144").Value = 2
$ws.Range("A145").Value = 2
$ws.Range("A147").Value = 1
$ws.Range("A148").Value = 2
$ws.Range("A151").Value = 3
$ws.Range("A152").Value = 2
$ws.Range("A153").Value = 1
$ws.Range("A156").Value = 3
$ws.Range("A157").Value = 1
$ws.Range("A160").Value = 3
$ws.Range("A162").Value = 2
$ws.Range("A163").Value = 1
$ws.Range("A164").Value = 1
$ws.Range("A166").Value = 3
$ws.Range("A168").Value = 3
$ws.Range("A170").Value = 3
$ws.Range("A171").Value = 3
$ws.Range("A172").Value = 1
$ws.Range("A173").Value = 3
$ws.Range("A175").Value = 1
$ws.Range("A177").Value = 1
$ws.Range("A178").Value = 3
$ws.Range("A179").Value = 3
$ws.Range("A180").Value = 3
$ws.Range("A181").Value = 2
$ws.Range("A183").Value = 3
$ws.Range("A184").Value = 3
$ws.Range("A185").Value = 2
$ws.Range("A187").Value = 3
$ws.Range("A189").Value = 1
$ws.Range("A190").Value = 2
$ws.Range("A191").Value = 3
$ws.Range("A192").Value = 2
$ws.Range("A193").Value = 3
$ws.Range("A195").Value = 1
$ws.Range("A196").Value = 3
$ws.Range("A197").Value = 1
$ws.Range("A198").Value = 3
$ws.Range("A200").Value = 3

# Remove now-unused rows 202-251 (data trimmed to 250 rows total, 200 data rows)
$ws.Range("A202:A251").ClearContents() | Out-Null

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
